$d = $word.ActiveDocument

# Locate the range spanning the three runs "<id>", "p126r_a1", "</id>"
# (together they read "<id>p126r_a1</id>" in the document text).
$rng = $d.Content
$found = $rng.Find.Execute("<id>p126r_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text '<id>p126r_a1</id>' in document."
}

$start = $rng.Start
$end = $rng.End

# Replacing the Range's text collapses the three differently-formatted runs
# into a single run that inherits the formatting of the first run
# (Courier New / color 7f6000 / sz 18), matching the target markup
# "<id>p126r_1</id>".
$rng.Text = "<id>p126r_1</id>"

$check = $d.Range($start, $start + ("<id>p126r_1</id>".Length))
Write-Host "Updated range text: $($check.Text)"
